# Add a new "Save" column (H) to the s_vals sheet, matching the style
# of the existing header row and the plain (unstyled) numeric cells below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last header cell (G1) onto the new header
# cell (H1) so it picks up the same bold/centered/bordered style, then
# set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the new data column with the "Save" flag values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
